$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 772
$ws.Range("F3").Value = 2800
$ws.Range("F4").Value = 1331
$ws.Range("F5").Value = 62
$ws.Range("F8").Value = 48
$ws.Range("F9").Value = 605
$ws.Range("F10").Value = 282
$ws.Range("F12").Value = 11663
$ws.Range("F13").Value = 6634
$ws.Range("F14").Value = 24
$ws.Range("F16").Value = 419
$ws.Range("F17").Value = 254
$ws.Range("F18").Value = 9
$ws.Range("F20").Value = 920
$ws.Range("F21").Value = 83
$ws.Range("F22").Value = 269
$ws.Range("F23").Value = 926
$ws.Range("F24").Value = 3649
$ws.Range("F25").Value = 58
$ws.Range("F26").Value = 987
$ws.Range("F27").Value = 496
$ws.Range("F28").Value = 168
$ws.Range("F29").Value = 316
$ws.Range("F30").Value = 2
$ws.Range("F31").Value = 225
$ws.Range("F32").Value = 269
$ws.Range("F33").Value = 305
$ws.Range("F34").Value = 5017
$ws.Range("F35").Value = 42
$ws.Range("F36").Value = 1241
$ws.Range("F38").Value = 447
$ws.Range("F39").Value = 201
$ws.Range("F40").Value = 544

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 3677

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9051
$ws.Range("F3").Value = 501
$ws.Range("F4").Value = 1831

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 501
$ws.Range("F3").Value = 1831
$ws.Range("F4").Value = 772
$ws.Range("F5").Value = 2800
$ws.Range("F9").Value = 1331
$ws.Range("F11").Value = 62
$ws.Range("F15").Value = 605
$ws.Range("F16").Value = 282
$ws.Range("F18").Value = 11663
$ws.Range("F19").Value = 3677
$ws.Range("F20").Value = 6634
$ws.Range("F22").Value = 24
$ws.Range("F24").Value = 419
$ws.Range("F25").Value = 254
$ws.Range("F26").Value = 9
$ws.Range("F28").Value = 83
$ws.Range("F29").Value = 269
$ws.Range("F30").Value = 926
$ws.Range("F31").Value = 3649
$ws.Range("F32").Value = 58
$ws.Range("F33").Value = 987
$ws.Range("F34").Value = 168
$ws.Range("F35").Value = 316
$ws.Range("F36").Value = 225
$ws.Range("F37").Value = 269
$ws.Range("F40").Value = 42
$ws.Range("F41").Value = 1241
$ws.Range("F44").Value = 201
$ws.Range("F45").Value = 544
